$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.600.79'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.894.01'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.99'
$ws.Range('E5').Value = '  +3.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4930'
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2950'
$ws.Range('E8').Value = '  +2.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06725'
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').Value = '1.910.21'
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.24'
$ws.Range('E11').Value = '  +3.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07363'
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.166'
$ws.Range('E13').Value = '  +4.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.38'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6724'
$ws.Range('E15').Value = '  +2.18%  '
$ws.Range('D16').Value = '30.538.75'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007884'
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('E18').Value = '  +4.89%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = '2.131.93'
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.330'
$ws.Range('E21').Value = '  +13.26%  '
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '191.28'
$ws.Range('E23').Value = '  +3.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.224'
$ws.Range('E24').Value = '  +3.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.557'
$ws.Range('E25').Value = '  +3.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.06'
$ws.Range('E26').Value = '  +3.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.45'
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.945'
$ws.Range('E28').Value = '  +6.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.474'
$ws.Range('E29').Value = '  +5.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.452'
$ws.Range('E30').Value = '  +5.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09194'
$ws.Range('E31').Value = '  +2.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.206'
$ws.Range('E32').Value = '  +8.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05251'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7442'
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.105'
$ws.Range('E35').Value = '  +3.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.715'
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01835'
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.693'
$ws.Range('E38').Value = '  +1.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9228'
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.069'
$ws.Range('E40').Value = '  +3.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4440'
$ws.Range('E41').Value = '  +3.81%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.955'
$ws.Range('E42').Value = '  +6.48%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.02'
$ws.Range('E43').Value = '  +27.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '106.23'
$ws.Range('E44').Value = '  +2.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9934'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1389'
$ws.Range('E46').Value = '  +4.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.596'
$ws.Range('E47').Value = '  +5.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.028'
$ws.Range('E48').Value = '  +5.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.14'
$ws.Range('E49').Value = '  +6.63%  '
$ws.Range('E50').Value = '  +0.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3975'
$ws.Range('E51').Value = '  +3.26%  '
